$d = $word.ActiveDocument

# Locate the end of the paragraph that ends with the "sidelines" example
# sentence; the five new paragraphs must be inserted right after it (and
# right before the existing blank separator paragraph that follows it).
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "She could only watch from the sidelines as her brother’s health deteriorated.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor sentence in document"
}

$insertPos = $findRange.End
$insertionPoint = $d.Range($insertPos, $insertPos)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# New blank separator paragraph (indented like the examples above it).
$para1 = "<w:p xmlns:w='$wNs'><w:pPr><w:ind w:left=`"360`"/></w:pPr></w:p>"

# New bulleted vocabulary heading + definition paragraph.
$para2 = "<w:p xmlns:w='$wNs'>" +
         "<w:pPr>" +
         "<w:pStyle w:val=`"ListParagraph`"/>" +
         "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
         "</w:pPr>" +
         "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`">In favor of something/doing something: </w:t></w:r>" +
         "<w:r><w:t>in a way that supports something or helps it to be successful:</w:t></w:r>" +
         "</w:p>"

# Three example paragraphs, each with a Wingdings arrow symbol + sentence.
$para3 = "<w:p xmlns:w='$wNs'><w:pPr><w:ind w:left=`"360`"/></w:pPr>" +
         "<w:r><w:sym w:font=`"Wingdings`" w:char=`"F0E0`"/></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> The members voted in favor of resolution.</w:t></w:r>" +
         "</w:p>"

$para4 = "<w:p xmlns:w='$wNs'><w:pPr><w:ind w:left=`"360`"/></w:pPr>" +
         "<w:r><w:sym w:font=`"Wingdings`" w:char=`"F0E0`"/></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> She spoke in favor of increasing the allocation of funds to rural districts.</w:t></w:r>" +
         "</w:p>"

$para5 = "<w:p xmlns:w='$wNs'><w:pPr><w:ind w:left=`"360`"/></w:pPr>" +
         "<w:r><w:sym w:font=`"Wingdings`" w:char=`"F0E0`"/></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> He complained that the system was fixed in favor of the large corporations.</w:t></w:r>" +
         "</w:p>"

# All five fragments must be inserted together in one InsertXML call: the
# paragraph immediately following the anchor sentence is empty (no runs),
# and inserting a lone <w:p> fragment right at that boundary gets merged
# into the following empty paragraph instead of becoming its own
# paragraph. Batching the fragments avoids that edge case.
$xml = $para1 + $para2 + $para3 + $para4 + $para5

$insertionPoint.InsertXML($xml)

Write-Output "Inserted new vocabulary entry after sidelines example."
